$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,14
$row2[0,0] = 16.02346907382236
$row2[0,1] = 0
$row2[0,2] = 10.95531719642747
$row2[0,3] = 17.36141455613541
$row2[0,4] = 32.85759418847496
$row2[0,5] = 32.71920215212118
$row2[0,6] = 15.38509129843432
$row2[0,7] = 0
$row2[0,8] = 12.20774103681907
$row2[0,9] = 8.695547329368882
$row2[0,10] = 8.241448496448983
$row2[0,11] = 14.90628147891183
$row2[0,12] = 21.0548770987178
$row2[0,13] = 23.92663729342057
$ws.Range("B2:O2").Value = $row2

$row3 = New-Object 'object[,]' 1,14
$row3[0,0] = 15.94057154318303
$row3[0,1] = 0
$row3[0,2] = 10.96875836165673
$row3[0,3] = 17.40444044114252
$row3[0,4] = 32.93386660921899
$row3[0,5] = 32.78507376946699
$row3[0,6] = 15.42538691167375
$row3[0,7] = 0
$row3[0,8] = 12.22902463626174
$row3[0,9] = 8.416329389905506
$row3[0,10] = 8.215499193498864
$row3[0,11] = 14.88137961979276
$row3[0,12] = 21.1109693197901
$row3[0,13] = 23.99123736773323
$ws.Range("B3:O3").Value = $row3

$row4 = New-Object 'object[,]' 1,14
$row4[0,0] = 15.89210916243093
$row4[0,1] = 0
$row4[0,2] = 10.97834967378302
$row4[0,3] = 17.43249614580302
$row4[0,4] = 32.98674443285262
$row4[0,5] = 32.83345956070339
$row4[0,6] = 15.45211269182258
$row4[0,7] = 0
$row4[0,8] = 12.24278897442126
$row4[0,9] = 8.238516468244214
$row4[0,10] = 8.200297823639966
$row4[0,11] = 14.86779089872029
$row4[0,12] = 21.14710512716267
$row4[0,13] = 24.03495031608647
$ws.Range("B4:O4").Value = $row4

$row5 = New-Object 'object[,]' 1,14
$row5[0,0] = 15.872989244053
$row5[0,1] = 0
$row5[0,2] = 10.98259533587075
$row5[0,3] = 17.44434176638642
$row5[0,4] = 33.00981201436617
$row5[0,5] = 32.85517022022483
$row5[0,6] = 15.4635029815653
$row5[0,7] = 0
$row5[0,8] = 12.24857358609922
$row5[0,9] = 8.164521644710772
$row5[0,10] = 8.194290828911232
$row5[0,11] = 14.86268523267548
$row5[0,12] = 21.16225813344377
$row5[0,13] = 24.05378121486026
$ws.Range("B5:O5").Value = $row5

$row6 = New-Object 'object[,]' 1,14
$row6[0,0] = 15.86985283377998
$row6[0,1] = 0
$row6[0,2] = 10.98332070249145
$row6[0,3] = 17.44633367472531
$row6[0,4] = 33.01373410807165
$row6[0,5] = 32.85889551773839
$row6[0,6] = 15.46542449965912
$row6[0,7] = 0
$row6[0,8] = 12.24954473249664
$row6[0,9] = 8.152144284769635
$row6[0,10] = 8.193304813555407
$row6[0,11] = 14.8618636381146
$row6[0,12] = 21.16480012497393
$row6[0,13] = 24.05696951272176
$ws.Range("B6:O6").Value = $row6

$row7 = New-Object 'object[,]' 1,14
$row7[0,0] = 15.89184873751158
$row7[0,1] = 0
$row7[0,2] = 10.97840556651689
$row7[0,3] = 17.43265422780364
$row7[0,4] = 32.98704937954252
$row7[0,5] = 32.83374429363049
$row7[0,6] = 15.45226428301136
$row7[0,7] = 0
$row7[0,8] = 12.24286627634429
$row7[0,9] = 8.237524667383671
$row7[0,10] = 8.200216046469629
$row7[0,11] = 14.8677202881584
$row7[0,12] = 21.14730775394901
$row7[0,13] = 24.03520015691135
$ws.Range("B7:O7").Value = $row7

$row8 = New-Object 'object[,]' 1,14
$row8[0,0] = 15.99439008756526
$row8[0,1] = 0
$row8[0,2] = 10.95967428535277
$row8[0,3] = 17.3759105291
$row8[0,4] = 32.88263750761137
$row8[0,5] = 32.74026466244997
$row8[0,6] = 15.39857366373769
$row8[0,7] = 0
$row8[0,8] = 12.21493543533012
$row8[0,9] = 8.600634000427373
$row8[0,10] = 8.232351604151596
$row8[0,11] = 14.8973447678994
$row8[0,12] = 21.07386660933665
$row8[0,13] = 23.94807093519616
$ws.Range("B8:O8").Value = $row8

$row9 = New-Object 'object[,]' 1,14
$row9[0,0] = 16.21406816514233
$row9[0,1] = 0
$row9[0,2] = 10.93353555256661
$row9[0,3] = 17.27758964607088
$row9[0,4] = 32.72589688318653
$row9[0,5] = 32.62009325224378
$row9[0,6] = 15.30901178976371
$row9[0,7] = 0
$row9[0,8] = 12.16566426423333
$row9[0,9] = 9.259496940457938
$row9[0,10] = 8.301005715987928
$row9[0,11] = 14.96874819267278
$row9[0,12] = 20.9432438158352
$row9[0,13] = 23.80935028149906
$ws.Range("B9:O9").Value = $row9

$row10 = New-Object 'object[,]' 1,14
$row10[0,0] = 16.38576059111344
$row10[0,1] = 0
$row10[0,2] = 10.92075283993218
$row10[0,3] = 17.21319363063759
$row10[0,4] = 32.64004953890168
$row10[0,5] = 32.57045182493779
$row10[0,6] = 15.25277282049056
$row10[0,7] = 0
$row10[0,8] = 12.13278703542827
$row10[0,9] = 9.70805331768481
$row10[0,10] = 8.354648766063153
$row10[0,11] = 15.02905894881387
$row10[0,12] = 20.85536515964304
$row10[0,13] = 23.72704835238798
$ws.Range("B10:O10").Value = $row10

$row11 = New-Object 'object[,]' 1,14
$row11[0,0] = 16.46586882422807
$row11[0,1] = 0
$row11[0,2] = 10.91632344295565
$row11[0,3] = 17.18558876786279
$row11[0,4] = 32.6073645187871
$row11[0,5] = 32.55628079206939
$row11[0,6] = 15.22925939133684
$row11[0,7] = 0
$row11[0,8] = 12.11854507939742
$row11[0,9] = 9.903867654680463
$row11[0,10] = 8.379695537329344
$row11[0,11] = 15.05813985836389
$row11[0,12] = 20.81712735329164
$row11[0,13] = 23.69387028426747
$ws.Range("B11:O11").Value = $row11

$row12 = New-Object 'object[,]' 1,14
$row12[0,0] = 16.49647142257204
$row12[0,1] = 0
$row12[0,2] = 10.91484458133261
$row12[0,3] = 17.17537753614275
$row12[0,4] = 32.59590316193403
$row12[0,5] = 32.55212451298203
$row12[0,6] = 15.22065273269667
$row12[0,7] = 0
$row12[0,8] = 12.11325420347826
$row12[0,9] = 9.976795322891725
$row12[0,10] = 8.389268104690489
$row12[0,11] = 15.06938308937035
$row12[0,12] = 20.8028965143047
$row12[0,13] = 23.68191962121626
$ws.Range("B12:O12").Value = $row12

$row13 = New-Object 'object[,]' 1,14
$row13[0,0] = 16.48986905186797
$row13[0,1] = 0
$row13[0,2] = 10.91515426725191
$row13[0,3] = 17.17756595102301
$row13[0,4] = 32.59833084246052
$row13[0,5] = 32.55296582791607
$row13[0,6] = 15.22249310984044
$row13[0,7] = 0
$row13[0,8] = 12.1143891478754
$row13[0,9] = 9.961144001274013
$row13[0,10] = 8.387202648434446
$row13[0,11] = 15.06695148051108
$row13[0,12] = 20.80595032247434
$row13[0,13] = 23.68446613548235
$ws.Range("B13:O13").Value = $row13

$row14 = New-Object 'object[,]' 1,14
$row14[0,0] = 16.46838128140795
$row14[0,1] = 0
$row14[0,2] = 10.91619780342979
$row14[0,3] = 17.18474383543263
$row14[0,4] = 32.60640323316749
$row14[0,5] = 32.55591460248612
$row14[0,6] = 15.22854535836003
$row14[0,7] = 0
$row14[0,8] = 12.11810774922413
$row14[0,9] = 9.909892172340477
$row14[0,10] = 8.380481344157001
$row14[0,11] = 15.05906025091609
$row14[0,12] = 20.81595158933418
$row14[0,13] = 23.6928748058183
$ws.Range("B14:O14").Value = $row14

$row15 = New-Object 'object[,]' 1,14
$row15[0,0] = 16.45525359144063
$row15[0,1] = 0
$row15[0,2] = 10.9168628201666
$row15[0,3] = 17.18917200659821
$row15[0,4] = 32.61146706307925
$row15[0,5] = 32.55787838929302
$row15[0,6] = 15.2322912525512
$row15[0,7] = 0
$row15[0,8] = 12.12039880216958
$row15[0,9] = 9.878338615402214
$row15[0,10] = 8.376375655496197
$row15[0,11] = 15.05425654608491
$row15[0,12] = 20.82211004948674
$row15[0,13] = 23.6981052187301
$ws.Range("B15:O15").Value = $row15

$row16 = New-Object 'object[,]' 1,14
$row16[0,0] = 16.38056388835792
$row16[0,1] = 0
$row16[0,2] = 10.921070125947
$row16[0,3] = 17.21503159804152
$row16[0,4] = 32.64231372633123
$row16[0,5] = 32.57154724561327
$row16[0,6] = 15.25435110647342
$row16[0,7] = 0
$row16[0,8] = 12.13373211199405
$row16[0,9] = 9.695087246039339
$row16[0,10] = 8.353024461795609
$row16[0,11] = 15.02719107034351
$row16[0,12] = 20.85789898808529
$row16[0,13] = 23.72930240135379
$ws.Range("B16:O16").Value = $row16

$row17 = New-Object 'object[,]' 1,14
$row17[0,0] = 16.33524283900412
$row17[0,1] = 0
$row17[0,2] = 10.92400549597116
$row17[0,3] = 17.23132771195243
$row17[0,4] = 32.66286816116692
$row17[0,5] = 32.58208750204887
$row17[0,6] = 15.26841407278757
$row17[0,7] = 0
$row17[0,8] = 12.14209424122929
$row17[0,9] = 9.580529961772086
$row17[0,10] = 8.33886093261709
$row17[0,11] = 15.01100470355719
$row17[0,12] = 20.88029891558946
$row17[0,13] = 23.74953260972945
$ws.Range("B17:O17").Value = $row17

$row18 = New-Object 'object[,]' 1,14
$row18[0,0] = 16.30936532097521
$row18[0,1] = 0
$row18[0,2] = 10.92582430485263
$row18[0,3] = 17.24085984769282
$row18[0,4] = 32.67528986259555
$row18[0,5] = 32.58894168874856
$row18[0,6] = 15.27669757110407
$row18[0,7] = 0
$row18[0,8] = 12.14697115341153
$row18[0,9] = 9.513866398814965
$row18[0,10] = 8.330775379059622
$row18[0,11] = 15.00184994430114
$row18[0,12] = 20.89334647052128
$row18[0,13] = 23.76156958609991
$ws.Range("B18:O18").Value = $row18

$row19 = New-Object 'object[,]' 1,14
$row19[0,0] = 16.3006368920973
$row19[0,1] = 0
$row19[0,2] = 10.92646255226157
$row19[0,3] = 17.24411460445469
$row19[0,4] = 32.67959856550398
$row19[0,5] = 32.59139834906381
$row19[0,6] = 15.27953569614498
$row19[0,7] = 0
$row19[0,8] = 12.14863395342327
$row19[0,9] = 9.491163712304912
$row19[0,10] = 8.328048357881984
$row19[0,11] = 14.9987771267414
$row19[0,12] = 20.89779229706643
$row19[0,13] = 23.76571397595325
$ws.Range("B19:O19").Value = $row19

$row20 = New-Object 'object[,]' 1,14
$row20[0,0] = 16.34004783294044
$row20[0,1] = 0
$row20[0,2] = 10.92367952404959
$row20[0,3] = 17.22957650820437
$row20[0,4] = 32.66061807593148
$row20[0,5] = 32.5808835289676
$row20[0,6] = 15.26689688165633
$row20[0,7] = 0
$row20[0,8] = 12.14119712274218
$row20[0,9] = 9.592805102869319
$row20[0,10] = 8.340362392108078
$row20[0,11] = 15.01271174653975
$row20[0,12] = 20.87789746925361
$row20[0,13] = 23.74733755480714
$ws.Range("B20:O20").Value = $row20

$row21 = New-Object 'object[,]' 1,14
$row21[0,0] = 16.47468567314179
$row21[0,1] = 0
$row21[0,2] = 10.91588591200976
$row21[0,3] = 17.18262895156488
$row21[0,4] = 32.60400732329278
$row21[0,5] = 32.55501563714181
$row21[0,6] = 15.22675959815291
$row21[0,7] = 0
$row21[0,8] = 12.11701273446659
$row21[0,9] = 9.92497955219496
$row21[0,10] = 8.382453206966483
$row21[0,11] = 15.06137187536134
$row21[0,12] = 20.81300722546012
$row21[0,13] = 23.69038833233897
$ws.Range("B21:O21").Value = $row21

$row22 = New-Object 'object[,]' 1,14
$row22[0,0] = 16.56422764958502
$row22[0,1] = 0
$row22[0,2] = 10.91194868747761
$row22[0,3] = 17.15335697889447
$row22[0,4] = 32.57234650739446
$row22[0,5] = 32.54516188520844
$row22[0,6] = 15.2022607059087
$row22[0,7] = 0
$row22[0,8] = 12.10180259346241
$row22[0,9] = 10.13493169221
$row22[0,10] = 8.410472053842261
$row22[0,11] = 15.09451729699904
$row22[0,12] = 20.77204852777834
$row22[0,13] = 23.65674274519604
$ws.Range("B22:O22").Value = $row22

$row23 = New-Object 'object[,]' 1,14
$row23[0,0] = 16.51630273580804
$row23[0,1] = 0
$row23[0,2] = 10.91394452229426
$row23[0,3] = 17.16885113166412
$row23[0,4] = 32.58875611173507
$row23[0,5] = 32.54977575439799
$row23[0,6] = 15.2151777402272
$row23[0,7] = 0
$row23[0,8] = 12.10986616376143
$row23[0,9] = 10.02354142530588
$row23[0,10] = 8.395472792118296
$row23[0,11] = 15.07670598230806
$row23[0,12] = 20.79377654033481
$row23[0,13] = 23.67437290896578
$ws.Range("B23:O23").Value = $row23

$row24 = New-Object 'object[,]' 1,14
$row24[0,0] = 16.3378749386612
$row24[0,1] = 0
$row24[0,2] = 10.92382648715564
$row24[0,3] = 17.23036771876123
$row24[0,4] = 32.66163345595194
$row24[0,5] = 32.58142537059523
$row24[0,6] = 15.26758218539514
$row24[0,7] = 0
$row24[0,8] = 12.14160249376597
$row24[0,9] = 9.587258009958067
$row24[0,10] = 8.339683403592787
$row24[0,11] = 15.01193952163833
$row24[0,12] = 20.87898263496332
$row24[0,13] = 23.74832867338642
$ws.Range("B24:O24").Value = $row24

$row25 = New-Object 'object[,]' 1,14
$row25[0,0] = 16.15275356296805
$row25[0,1] = 0
$row25[0,2] = 10.93947629456534
$row25[0,3] = 17.30280716558086
$row25[0,4] = 32.76315438211488
$row25[0,5] = 32.64582438912784
$row25[0,6] = 15.33155966705655
$row25[0,7] = 0
$row25[0,8] = 12.17840774243448
$row25[0,9] = 9.087293228025512
$row25[0,10] = 8.281853075694512
$row25[0,11] = 14.94803301831235
$row25[0,12] = 20.97715461742925
$row25[0,13] = 23.84343449986902
$ws.Range("B25:O25").Value = $row25
